$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diacritic "shadda" (bare) entry in row 9 is being removed. Every row
# below it moves its diacritic value (and the wrap-text formatting that
# travels with it) up by one, the newly-vacated last diacritic row (15) is
# cleared, and the now fully empty row 16 is deleted outright.
for ($r = 9; $r -le 14; $r++) {
    $srcCell = $ws.Cells.Item($r + 1, 2)
    $dstCell = $ws.Cells.Item($r, 2)

    $srcCell.Copy()
    $dstCell.PasteSpecial(-4163)  # xlPasteValues
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

$ws.Cells.Item(15, 2).Clear()
$ws.Rows("16").Delete()

# Restore the view state captured in the saved workbook.
$ws.Range("B8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 175
